$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1563.4333
$ws.Range("J17").Value = 1423.1428
$ws.Range("L17").Value = 4269.428400000001
$ws.Range("N17").Value = -4605.428400000001

# Row 97
$ws.Range("H97").Value = 1738.8182
$ws.Range("J97").Value = 1738.8182
$ws.Range("L97").Value = 5216.4546
$ws.Range("N97").Value = -6208.4546

# Row 121
$ws.Range("H121").Value = 3649.5
$ws.Range("J121").Value = 5000
$ws.Range("L121").Value = 15000
$ws.Range("N121").Value = -18494

# Row 129
$ws.Range("H129").Value = 1965.6
$ws.Range("I129").Value = 548.5
$ws.Range("K129").Value = 1645.5
$ws.Range("M129").Value = 3354.5

# Row 135
$ws.Range("H135").Value = 1372.2106
$ws.Range("I135").Value = 1237.5
$ws.Range("J135").Value = 2090.6667
$ws.Range("K135").Value = 11137.5
$ws.Range("L135").Value = 18816.0003
$ws.Range("M135").Value = -8602.5
$ws.Range("N135").Value = -23886.0003

# Row 137
$ws.Range("H137").Value = 1284923.1
$ws.Range("J137").Value = 3132.0908
$ws.Range("L137").Value = 9396.2724
$ws.Range("N137").Value = -14496.2724

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 49113.11
$ws.Range("I45").Value = 49113.11
$ws.Range("K45").Value = 49113.11
$ws.Range("M45").Value = -48736.11

# Row 61
$ws.Range("H61").Value = 3010.4
$ws.Range("I61").Value = 2470.8
$ws.Range("K61").Value = 2470.8
$ws.Range("M61").Value = -2258.8

# Row 74
$ws.Range("H74").Value = 160728
$ws.Range("I74").Value = 170166.1
$ws.Range("K74").Value = 170166.1
$ws.Range("M74").Value = -169292.1

# Row 77
$ws.Range("H77").Value = 160728
$ws.Range("I77").Value = 170166.1
$ws.Range("K77").Value = 850830.5
$ws.Range("M77").Value = -846462.5

# Row 88
$ws.Range("H88").Value = 2265.8484
$ws.Range("I88").Value = 1736.909
$ws.Range("J88").Value = 2530.318
$ws.Range("K88").Value = 1736.909
$ws.Range("L88").Value = 2530.318
$ws.Range("M88").Value = -1330.909
$ws.Range("N88").Value = -3342.318

# Row 91
$ws.Range("H91").Value = 2265.8484
$ws.Range("I91").Value = 1736.909
$ws.Range("J91").Value = 2530.318
$ws.Range("K91").Value = 1736.909
$ws.Range("L91").Value = 2530.318
$ws.Range("M91").Value = -332.9090000000001
$ws.Range("N91").Value = -5338.318

# Row 136
$ws.Range("H136").Value = 3010.4
$ws.Range("I136").Value = 2470.8
$ws.Range("K136").Value = 7412.400000000001
$ws.Range("M136").Value = -4862.400000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 13893962
$ws.Range("J20").Value = 2592.9285
$ws.Range("L20").Value = 2592.9285
$ws.Range("N20").Value = -3086.9285

# Row 75
$ws.Range("H75").Value = 14851.714
$ws.Range("I75").Value = 15327
$ws.Range("J75").Value = 12000
$ws.Range("K75").Value = 15327
$ws.Range("L75").Value = 12000
$ws.Range("M75").Value = -14391
$ws.Range("N75").Value = -13872

# Row 78
$ws.Range("H78").Value = 14851.714
$ws.Range("I78").Value = 15327
$ws.Range("J78").Value = 12000
$ws.Range("K78").Value = 45981
$ws.Range("L78").Value = 36000
$ws.Range("M78").Value = -41301
$ws.Range("N78").Value = -45360

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4635422.5
$ws.Range("J31").Value = 20842834
$ws.Range("L31").Value = 20842834
$ws.Range("N31").Value = -20843424

# Row 34
$ws.Range("H34").Value = 4635422.5
$ws.Range("J34").Value = 20842834
$ws.Range("L34").Value = 20842834
$ws.Range("N34").Value = -20843238

# Row 132
$ws.Range("H132").Value = 17546810
$ws.Range("I132").Value = 2349.7693
$ws.Range("J132").Value = 55559810
$ws.Range("K132").Value = 7049.3079
$ws.Range("L132").Value = 166679430
$ws.Range("M132").Value = -4519.3079
$ws.Range("N132").Value = -166684490

# Row 133
$ws.Range("H133").Value = 79599.7
$ws.Range("J133").Value = 79599.7
$ws.Range("L133").Value = 79599.7
$ws.Range("N133").Value = -84659.7

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 721.7143
$ws.Range("I5").Value = 577.7143
$ws.Range("J5").Value = 865.7143
$ws.Range("K5").Value = 1733.1429
$ws.Range("L5").Value = 2597.1429
$ws.Range("M5").Value = -1621.1429
$ws.Range("N5").Value = -2821.1429

# Row 56
$ws.Range("H56").Value = 14638.75
$ws.Range("I56").Value = 14638.75
$ws.Range("K56").Value = 14638.75
$ws.Range("M56").Value = -14108.75

# Row 135
$ws.Range("H135").Value = 721.7143
$ws.Range("I135").Value = 577.7143
$ws.Range("J135").Value = 865.7143
$ws.Range("K135").Value = 5199.428699999999
$ws.Range("L135").Value = 7791.428699999999
$ws.Range("M135").Value = -2664.428699999999
$ws.Range("N135").Value = -12861.4287

# Row 141
$ws.Range("H141").Value = 4910.3687
$ws.Range("I141").Value = 4370.4116
$ws.Range("J141").Value = 9500
$ws.Range("K141").Value = 13111.2348
$ws.Range("L141").Value = 28500
$ws.Range("M141").Value = -7931.234800000002
$ws.Range("N141").Value = -38860

$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 15009
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29849

# Row 114
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

# Row 122
$ws.Range("H122").Value = 4556.143
$ws.Range("I122").Value = 2186.125
$ws.Range("J122").Value = 7716.1665
$ws.Range("K122").Value = 6558.375
$ws.Range("L122").Value = 23148.4995
$ws.Range("M122").Value = -4108.375
$ws.Range("N122").Value = -28048.4995

# Row 132
$ws.Range("H132").Value = 3393.7273
$ws.Range("I132").Value = 3353.1
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 10059.3
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -7529.299999999999
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4165.05
$ws.Range("I7").Value = 3635.647
$ws.Range("K7").Value = 3635.647
$ws.Range("M7").Value = -3523.647

# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()

# Row 55
$ws.Range("H55").Value = 519.1875
$ws.Range("I55").Value = 333.33334
$ws.Range("J55").Value = 758.1429000000001
$ws.Range("K55").Value = 333.33334
$ws.Range("L55").Value = 758.1429000000001
$ws.Range("M55").Value = -160.33334
$ws.Range("N55").Value = -1104.1429

# Row 61
$ws.Range("H61").Value = 5313.7744
$ws.Range("I61").Value = 1690.2693
$ws.Range("J61").Value = 24156
$ws.Range("K61").Value = 1690.2693
$ws.Range("L61").Value = 24156
$ws.Range("M61").Value = -1488.2693
$ws.Range("N61").Value = -24560

# Row 113
$ws.Range("H113").Value = 5313.7744
$ws.Range("I113").Value = 1690.2693
$ws.Range("J113").Value = 24156
$ws.Range("K113").Value = 1690.2693
$ws.Range("L113").Value = 24156
$ws.Range("M113").Value = 479.7307000000001
$ws.Range("N113").Value = -28496

# Row 126
$ws.Range("H126").Value = 4165.05
$ws.Range("I126").Value = 3635.647
$ws.Range("K126").Value = 10906.941
$ws.Range("M126").Value = -8436.940999999999

# Row 135
$ws.Range("H135").Value = 105992
$ws.Range("J135").Value = 105992
$ws.Range("L135").Value = 105992
$ws.Range("N135").Value = -116132

# Row 140
$ws.Range("H140").Value = 149379.89
$ws.Range("J140").Value = 149379.89
$ws.Range("L140").Value = 149379.89
$ws.Range("N140").Value = -159739.89

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 96
$ws.Range("H96").Value = 2230.3333
$ws.Range("I96").Value = 1970
$ws.Range("J96").Value = 2751
$ws.Range("K96").Value = 1970
$ws.Range("L96").Value = 2751
$ws.Range("M96").Value = -597
$ws.Range("N96").Value = -5497

# Row 132
$ws.Range("H132").Value = 5881.1665
$ws.Range("I132").Value = 5498.75
$ws.Range("J132").Value = 6646
$ws.Range("K132").Value = 16496.25
$ws.Range("L132").Value = 19938
$ws.Range("M132").Value = -13966.25
$ws.Range("N132").Value = -24998

# Row 136
$ws.Range("H136").Value = 14221.413
$ws.Range("I136").Value = 14906.854
$ws.Range("J136").Value = 8600.799999999999
$ws.Range("K136").Value = 44720.562
$ws.Range("L136").Value = 25802.4
$ws.Range("M136").Value = -42170.562
$ws.Range("N136").Value = -30902.4
